$d = $word.ActiveDocument

# Locate the unique run of text " p ˄ q" (the end of the paragraph that
# reads "˥(p → q)  p ˄ q") so we can split it into three runs:
#   " p ˄ "  +  "˥"  +  "q"
# turning the visible text into " p ˄ ˥q" while keeping every run's
# formatting (font, color, highlight) identical to the original run.
$finder = $d.Content
$found = $finder.Find.Execute(" p ˄ q", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    # Re-anchor to a fresh Range over the same span: re-using the Range
    # object that Find just matched makes InsertXML *insert* rather than
    # *replace* the matched span, so build a brand-new Range(start, end).
    $rng = $d.Range($finder.Start, $finder.End)
    $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData>' +
           '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:body>' +
           '<w:p>' +
             '<w:r>' +
               '<w:rPr>' +
                 '<w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/>' +
                 '<w:color w:val="385623" w:themeColor="accent6" w:themeShade="80"/>' +
                 '<w:highlight w:val="yellow"/>' +
               '</w:rPr>' +
               '<w:t xml:space="preserve"> p ˄ </w:t>' +
             '</w:r>' +
             '<w:r>' +
               '<w:rPr>' +
                 '<w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/>' +
                 '<w:color w:val="385623" w:themeColor="accent6" w:themeShade="80"/>' +
                 '<w:highlight w:val="yellow"/>' +
               '</w:rPr>' +
               '<w:t>˥</w:t>' +
             '</w:r>' +
             '<w:r>' +
               '<w:rPr>' +
                 '<w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/>' +
                 '<w:color w:val="385623" w:themeColor="accent6" w:themeShade="80"/>' +
                 '<w:highlight w:val="yellow"/>' +
               '</w:rPr>' +
               '<w:t>q</w:t>' +
             '</w:r>' +
           '</w:p>' +
           '</w:body>' +
           '</w:document>' +
           '</pkg:xmlData>' +
           '</pkg:part>' +
           '</pkg:package>'

    $rng.InsertXML($xml)
}
